$d = $word.ActiveDocument

# --- 1. "...interact with crimes..." -> "...interact with a database of crimes..." ---
$findRange = $d.Range(0, 0)
$findRange.Find.Execute("interact with ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($findRange.End, $findRange.End)
$insertPoint.InsertBefore("a database of ")

# --- 2. Remove "(each is on its own website) " from the scraping sentence ---
$d.Content.Find.Execute(
    "We plan to scrape the individual reports (each is on its own website) to build our database.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We plan to scrape the individual reports to build our database.", 2) | Out-Null

# --- 3. Drop the old, now-redundant "_GoBack" bookmark paragraph -------------
#        (it was an otherwise-empty paragraph right after "List of tables
#        with keys declared:") and insert a fresh blank paragraph after
#        "List of tables with keys declared:" in its place.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "List of tables with keys declared:`r" -and $p.Range.Text.Length -gt 1) {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text -eq "`r") {
            $next.Range.Delete()
        }
        break
    }
}

$d.Content.Find.Execute(
    "List of tables with keys declared:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "List of tables with keys declared:^p", 2) | Out-Null

# --- 4. Re-add "_GoBack" spanning the whole document (start of doc to end) --
$wholeDoc = $d.Range(0, $d.Content.End)
$d.Bookmarks.Add("_GoBack", $wholeDoc)
